$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.461.64"
$ws.Range("E2").Value = "  +10.39%  "
$ws.Range("D3").Value = "3.257.58"
$ws.Range("E3").Value = "  +6.39%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'398.34"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'111.20"
$ws.Range("E6").Value = "  +9.49%  "
$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  +4.67%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  +7.23%  "
$ws.Range("D10").Value = "'39.54"
$ws.Range("E10").Value = "  +7.71%  "
$ws.Range("D11").Value = "'0.0948"
$ws.Range("E11").Value = "  +11.95%  "
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "3.774.94"
$ws.Range("E13").Value = "  +6.57%  "
$ws.Range("D14").Value = "'19.27"
$ws.Range("E14").Value = "  +5.60%  "
$ws.Range("E15").Value = "  +6.05%  "
$ws.Range("D16").Value = "3.260.74"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("E17").Value = "  +5.82%  "
$ws.Range("D18").Value = "'11.07"
$ws.Range("E18").Value = "  +4.45%  "
$ws.Range("D19").Value = "56.436.95"
$ws.Range("E19").Value = "  +10.39%  "
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("E21").Value = "  +8.58%  "
$ws.Range("D22").Value = "'13.08"
$ws.Range("E22").Value = "  +7.05%  "
$ws.Range("D23").Value = "'298.01"
$ws.Range("E23").Value = "  +12.89%  "
$ws.Range("D24").Value = "'75.24"
$ws.Range("E24").Value = "  +8.14%  "
$ws.Range("D25").Value = "'3.23"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").Value = "'8.13"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("D27").Value = "'28.25"
$ws.Range("E27").Value = "  +5.77%  "
$ws.Range("E28").Value = "  +5.10%  "
$ws.Range("D29").Value = "'7.36"
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").Value = "'0.169"
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +6.97%  "
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("D34").Value = "'37.88"
$ws.Range("E34").Value = "  +6.34%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("D37").Value = "'51.86"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("D38").Value = "'3.12"
$ws.Range("E38").Value = "  +25.78%  "
$ws.Range("D39").Value = "'3.52"
$ws.Range("E39").Value = "  +5.39%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'17.60"
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("E42").Value = "  +6.40%  "
$ws.Range("D43").Value = "'133.57"
$ws.Range("E43").Value = "  +3.08%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.99"
$ws.Range("E44").Value = "  +6.15%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.120"
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("D46").Value = "'0.286"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'22.30"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "'2.17"
$ws.Range("E48").Value = "  +53.19%  "
$ws.Range("D49").Value = "2.149.52"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("D50").Value = "'2.10"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("E51").Value = "  -3.02%  "
